$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.115.73'
$ws.Range("E2").Value = '  +1.98%  '

$ws.Range("D3").Value = '3.824.71'
$ws.Range("E3").Value = '  +0.46%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '627.77'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.37'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").Value = '3.819.77'
$ws.Range("E7").Value = '  +0.40%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.454'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.43%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.67'
$ws.Range("D12").ClearFormats()

$ws.Range("E13").Value = '  +0.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.10'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.37%  '

$ws.Range("D15").Value = '4.469.75'
$ws.Range("E15").Value = '  +0.46%  '

$ws.Range("D16").Value = '3.793.01'
$ws.Range("E16").Value = '  -0.74%  '

$ws.Range("D17").Value = '69.134.34'
$ws.Range("E17").Value = '  +1.93%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.29'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.02%  '

$ws.Range("E19").Value = '  +1.26%  '

$ws.Range("E20").Value = '  +0.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '467.42'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.76%  '

$ws.Range("E22").Value = '  -1.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.711'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.20%  '

$ws.Range("E24").Value = '  +4.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.15'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.09'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.17'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.10'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.62%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("D30").Value = '3.979.36'
$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("E31").Value = '  +2.44%  '

$ws.Range("E32").Value = '  -3.83%  '

$ws.Range("E33").Value = '  -0.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.23'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.40%  '

$ws.Range("E35").Value = '  +0.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("E37").Value = '  +2.16%  '

$ws.Range("E38").Value = '  +7.73%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.34'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.92'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.13%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.979'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '156.16'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.21%  '

$ws.Range("E45").Value = '  +0.67%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.43'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.75%  '

$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.01'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.88%  '

$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.86'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.94%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.49'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.90'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '381.60'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.37%  '
